$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All target cells hold text (degrees/percentages/fractions) rather than
# numeric values, so force text format before assigning to avoid Excel
# auto-converting "30%" into a numeric percentage, "8 de 11" into a date, etc.

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

Set-TextValue $ws.Range("B2") "35°"
Set-TextValue $ws.Range("D2") "30%"

Set-TextValue $ws.Range("C3") "20°"
Set-TextValue $ws.Range("D3") "30%"
Set-TextValue $ws.Range("E3") "52%"

Set-TextValue $ws.Range("B4") "28°"
Set-TextValue $ws.Range("C4") "15°"
Set-TextValue $ws.Range("D4") "62%"
Set-TextValue $ws.Range("E4") "83%"

Set-TextValue $ws.Range("D5") "70%"
Set-TextValue $ws.Range("E5") "84%"
Set-TextValue $ws.Range("F5") "8 de 11"

Set-TextValue $ws.Range("D6") "66%"

Set-TextValue $ws.Range("E7") "68%"

Set-TextValue $ws.Range("C8") "22°"
Set-TextValue $ws.Range("D8") "35%"
Set-TextValue $ws.Range("E8") "49%"

Set-TextValue $ws.Range("B9") "33°"
Set-TextValue $ws.Range("D9") "37%"

Set-TextValue $ws.Range("B10") "31°"
Set-TextValue $ws.Range("D10") "49%"
Set-TextValue $ws.Range("E10") "74%"

Set-TextValue $ws.Range("B11") "31°"
Set-TextValue $ws.Range("D11") "50%"
